# Update countries & provincias Spain
# - Refresh the COVID case/recovery/death snapshot for a number of countries.
# - Swap the "Islas Malvinas" / "Montserrat" rows so Islas Malvinas sorts first.
# - Bump the "Datos actualizados a ..." timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Updated case / recovered / active / deaths counts for several countries ---

# Row 4  - Estados Unidos
$ws.Range("B4").Value = 8651581
$ws.Range("C4").Value = 66231
$ws.Range("D4").Value = 5646389
$ws.Range("E4").Value = 2776923
$ws.Range("G4").Value = 861
$ws.Range("H4").Value = 228269

# Row 11 - Colombia
$ws.Range("B11").Value = 990270
$ws.Range("C11").Value = 8570
$ws.Range("D11").Value = 893712
$ws.Range("E11").Value = 66922
$ws.Range("G11").Value = 172
$ws.Range("H11").Value = 29636

# Row 12 - Peru
$ws.Range("B12").Value = 879876
$ws.Range("C12").Value = 2991
$ws.Range("D12").Value = 796719
$ws.Range("E12").Value = 49173
$ws.Range("G12").Value = 47
$ws.Range("H12").Value = 33984

# Row 20 - Alemania
$ws.Range("B20").Value = 403874
$ws.Range("C20").Value = 12519
$ws.Range("D20").Value = 306100
$ws.Range("E20").Value = 87730

# Row 33 - Canada
$ws.Range("B33").Value = 209148
$ws.Range("C33").Value = 3194
$ws.Range("D33").Value = 175805
$ws.Range("E33").Value = 23481

# Row 48 - Egipto
$ws.Range("B48").Value = 106060
$ws.Range("C48").Value = 177
$ws.Range("D48").Value = 98624
$ws.Range("E48").Value = 1270
$ws.Range("G48").Value = 11
$ws.Range("H48").Value = 6166

# Row 64 - Nigeria
$ws.Range("B64").Value = 61805
$ws.Range("C64").Value = 138
$ws.Range("D64").Value = 56985
$ws.Range("E64").Value = 3693
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 1127

# Row 84 - Bulgaria
$ws.Range("B84").Value = 34930
$ws.Range("C84").Value = 1595
$ws.Range("D84").Value = 17833
$ws.Range("E84").Value = 16033
$ws.Range("G84").Value = 16
$ws.Range("H84").Value = 1064

# Row 97 - Noruega
$ws.Range("B97").Value = 17234
$ws.Range("C97").Value = 270
$ws.Range("E97").Value = 5092

# Row 113 - Haiti
$ws.Range("B113").Value = 9007
$ws.Range("C113").Value = 28
$ws.Range("D113").Value = 7311
$ws.Range("E113").Value = 1465

# Row 150 - Gambia
$ws.Range("B150").Value = 3659
$ws.Range("C150").Value = 2
$ws.Range("D150").Value = 2660
$ws.Range("E150").Value = 880
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 119

# Row 167 - Niger
$ws.Range("B167").Value = 1215
$ws.Range("C167").Value = 1
$ws.Range("E167").Value = 18

# --- Reorder "Islas Malvinas" ahead of "Montserrat" (rows 216/217) ---
# Capture the two full rows first, then swap them so Islas Malvinas'
# row ends up first, matching the new shared-string order.
$name216 = $ws.Range("A216").Value()
$b216 = $ws.Range("B216").Value()
$c216 = $ws.Range("C216").Value()
$d216 = $ws.Range("D216").Value()
$e216 = $ws.Range("E216").Value()
$f216 = $ws.Range("F216").Value()
$g216 = $ws.Range("G216").Value()
$h216 = $ws.Range("H216").Value()

$name217 = $ws.Range("A217").Value()
$b217 = $ws.Range("B217").Value()
$c217 = $ws.Range("C217").Value()
$d217 = $ws.Range("D217").Value()
$e217 = $ws.Range("E217").Value()
$f217 = $ws.Range("F217").Value()
$g217 = $ws.Range("G217").Value()
$h217 = $ws.Range("H217").Value()

$ws.Range("A216").Value = $name217
$ws.Range("B216").Value = $b217
$ws.Range("C216").Value = $c217
$ws.Range("D216").Value = $d217
$ws.Range("E216").Value = $e217
$ws.Range("F216").Value = $f217
$ws.Range("G216").Value = $g217
$ws.Range("H216").Value = $h217

$ws.Range("A217").Value = $name216
$ws.Range("B217").Value = $b216
$ws.Range("C217").Value = $c216
$ws.Range("D217").Value = $d216
$ws.Range("E217").Value = $e216
$ws.Range("F217").Value = $f216
$ws.Range("G217").Value = $g216
$ws.Range("H217").Value = $h216

# --- Bump the "last updated" timestamp text (cell A1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 00:32"
